# Update the IL NAICS emissions percentages sheet:
# - Column A currently holds "sector" values (Other/Food and Beverage/Chemicals).
#   Move those values to a new column E (headed "sector"), and repurpose
#   column A to hold NAICS codes (headed "naics_code").
# - Row 2 ("Other") has no corresponding NAICS code, so A2 stays blank.
# - Fix a mislabeled description in row 17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Capture current sector values (column A, rows 2-19) before overwrite.
$sectors = @{}
for ($r = 2; $r -le 19; $r++) {
    $sectors[$r] = $ws.Cells.Item($r, 1).Value()
}

# Step 2: Write header row for new layout (A..E)
$ws.Cells.Item(1, 1).Value = "naics_code"
$ws.Cells.Item(1, 2).Value = "description"
$ws.Cells.Item(1, 3).Value = "co2e_total"
$ws.Cells.Item(1, 4).Value = "percent_of_total"
$ws.Cells.Item(1, 5).Value = "sector"

# Step 3: naics codes for each row (2-19); row 2 ("Other") has no code.
$naicsCodes = @{
    2  = $null
    3  = "311611"
    4  = "312140"
    5  = "311225"
    6  = "311224"
    7  = "311221"
    8  = "325199"
    9  = "325998"
    10 = "325414"
    11 = "325194"
    12 = "325193"
    13 = "325120"
    14 = "325411"
    15 = "325311"
    16 = "325180"
    17 = "325613"
    18 = "325110"
    19 = "325211"
}

for ($r = 2; $r -le 19; $r++) {
    $code = $naicsCodes[$r]
    if ($null -ne $code) {
        $ws.Cells.Item($r, 1).Value = "'" + $code
    }
    $ws.Cells.Item($r, 5).Value = $sectors[$r]
}
$ws.Cells.Item(2, 1).ClearContents()

# Step 4: Fix the description text for row 17
$ws.Cells.Item(17, 2).Value = "Other Chemicals Manufacturing"
